$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of column J ("SW(S*)") across the 10 instances -------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11
$ws.Range("J12").Font.Name = "Calibri"

# --- Rows 14-17: summary labels (col A) + stats (col B) -------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting for the summary block: bold, size 12, vertically centred
$summary = $ws.Range("B14:B17")
$summary.Font.Bold = $true
$summary.Font.Size = 12
$summary.VerticalAlignment = -4108

$ws.Range("A14:B17").RowHeight = 15.6

# Page setup (portrait, A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select
